# 912 work by Alex
# Adds three new "Race" columns to the end of the header row on Sheet1:
#   CV1 -> "Race Unknown"
#   CW1 -> "Race Other"
#   CX1 -> "Race Refused to Answer"
# and moves the active selection to CX6 (matching the post-edit cursor
# position recorded in the workbook view).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("CV1").Value = "Race Unknown"
$ws.Range("CW1").Value = "Race Other"
$ws.Range("CX1").Value = "Race Refused to Answer"

# Move/record the selection like the author's final cursor position.
$ws.Range("CX6").Select()

# The author also toggled the workbook's formula reference style to R1C1
# (xl/workbook.xml calcPr/@refMode="R1C1").
$excel.ReferenceStyle = 4150
